$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Лист1" -> "тарифы"): fix the data plan text in E2 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E2").Value = "50ГБ, 1200 минут по РФ, 500 SMS\n Безлимитные соцсети и мессенджеры"
$ws1.Range("C21").Select()
$ws1.Name = "тарифы"

# --- Sheet3 ("Лист3") is no longer used: delete it ---
$wb.Worksheets.Item(3).Delete()

# --- Sheet2 ("Лист2" -> "пользователи"): build the new users table ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "пользователи"

$ws2.Range("A1").Value = "id"
$ws2.Range("A1").NumberFormat = "0.00"

$ws2.Range("B1").Value = "имя"
$ws2.Range("B1").NumberFormat = "@"

$ws2.Range("C1").Value = "Интернет + ТВ + Мобильная связь"
$ws2.Range("C1").NumberFormat = "0.00"
$ws2.Range("C1").Font.ThemeColor = 1
$ws2.Range("C1").VerticalAlignment = -4108

$ws2.Range("D1").Value = "Интернет"
$ws2.Range("D1").NumberFormat = "0.00"
$ws2.Range("D1").Font.ThemeColor = 1
$ws2.Range("D1").VerticalAlignment = -4108

$ws2.Range("E1").Value = "Интернет + ТВ"
$ws2.Range("E1").NumberFormat = "0.00"
$ws2.Range("E1").Font.ThemeColor = 1
$ws2.Range("E1").VerticalAlignment = -4108

$ws2.Columns.Item(1).ColumnWidth = 20.7109375
$ws2.Columns.Item(2).ColumnWidth = 25.7109375
$ws2.Columns.Item(3).ColumnWidth = 31.7109375
$ws2.Columns.Item(4).ColumnWidth = 10.7109375
$ws2.Columns.Item(5).ColumnWidth = 18.7109375

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws2.Range("A2").Select()
$ws2.Activate()
